# Insert 3 new weekly records at the top of the Uva price block (rows 358-360),
# pushing the existing rows 358-454 down to 361-457.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("358:360").Insert()

# New row 358: Flame Seedless, Provincia de Copiapó, $/bandeja 10 kilos
$ws.Range("A358").Value = 3
$ws.Range("B358").Value = "Femacal de La Calera"
$ws.Range("C358").Value = "Coquimbo"
$ws.Range("D358").Value = 44551
$ws.Range("E358").Value = 5
$ws.Range("F358").Value = "Fruta"
$ws.Range("G358").Value = 100109
$ws.Range("H358").Value = "Uva"
$ws.Range("I358").Value = 100109001
$ws.Range("J358").Value = "Uva"
$ws.Range("K358").Value = "Flame Seedless"
$ws.Range("L358").Value = "Primera"
$ws.Range("M358").Value = 50
$ws.Range("N358").Value = 10000
$ws.Range("O358").Value = 10000
$ws.Range("P358").Value = 10000
$ws.Range("Q358").Value = "`$/bandeja 10 kilos"
$ws.Range("R358").Value = "Provincia de Copiapó"
$ws.Range("S358").Value = 1000
$ws.Range("T358").Value = 10

# New row 359: Red Globe, Provincia de Copiapó, $/bandeja 10 kilos
$ws.Range("A359").Value = 3
$ws.Range("B359").Value = "Femacal de La Calera"
$ws.Range("C359").Value = "Coquimbo"
$ws.Range("D359").Value = 44551
$ws.Range("E359").Value = 5
$ws.Range("F359").Value = "Fruta"
$ws.Range("G359").Value = 100109
$ws.Range("H359").Value = "Uva"
$ws.Range("I359").Value = 100109001
$ws.Range("J359").Value = "Uva"
$ws.Range("K359").Value = "Red Globe"
$ws.Range("L359").Value = "Primera"
$ws.Range("M359").Value = 67
$ws.Range("N359").Value = 16000
$ws.Range("O359").Value = 16000
$ws.Range("P359").Value = 16000
$ws.Range("Q359").Value = "`$/bandeja 10 kilos"
$ws.Range("R359").Value = "Provincia de Copiapó"
$ws.Range("S359").Value = 1600
$ws.Range("T359").Value = 10

# New row 360: Superior Seedless, Provincia de Copiapó, $/bandeja 10 kilos
$ws.Range("A360").Value = 3
$ws.Range("B360").Value = "Femacal de La Calera"
$ws.Range("C360").Value = "Coquimbo"
$ws.Range("D360").Value = 44551
$ws.Range("E360").Value = 5
$ws.Range("F360").Value = "Fruta"
$ws.Range("G360").Value = 100109
$ws.Range("H360").Value = "Uva"
$ws.Range("I360").Value = 100109001
$ws.Range("J360").Value = "Uva"
$ws.Range("K360").Value = "Superior Seedless"
$ws.Range("L360").Value = "Primera"
$ws.Range("M360").Value = 65
$ws.Range("N360").Value = 16000
$ws.Range("O360").Value = 16000
$ws.Range("P360").Value = 16000
$ws.Range("Q360").Value = "`$/bandeja 10 kilos"
$ws.Range("R360").Value = "Provincia de Copiapó"
$ws.Range("S360").Value = 1600
$ws.Range("T360").Value = 10
